$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Three more result blocks appended below the existing FY1/FY2/FY3/total rows
# (rows 7-10). Each block is 4 rows: FY1, FY2, FY3, total - same layout as
# rows 7-10, starting at row 11.
$newRows = @(
    @("FY1", 5.3, 140, 0.879, 0.005, 17922.58, 11.372, 1800, 17923.277, 11.437, 1800),
    @("FY2", 204.45, 138.335, 13.554, 9.970000000000001, 10293.109, 623.92, 1400, 9037.578, 53.062, 912.592),
    @("FY3", 199.812, 90.834, 0.333, 0.005, 5971.546, -3010.251, 700, 5971.119, -3010.405, 700),
    @("总有效遮蔽并集时长(s)", 9.138),

    @("FY1", 5.3, 140, 0.879, 0.005, 17922.58, 11.372, 1800, 17923.277, 11.437, 1800),
    @("FY2", 204.45, 138.335, 13.554, 9.970000000000001, 10293.109, 623.92, 1400, 9037.578, 53.062, 912.592),
    @("FY3", 199.812, 90.834, 0.333, 0.005, 5971.546, -3010.251, 700, 5971.119, -3010.405, 700),
    @("总有效遮蔽并集时长(s)", 9.138),

    @("FY1", 5.326, 140, 0.881, 0.005, 17922.874, 11.454, 1800, 17923.571, 11.519, 1800),
    @("FY2", 204.53, 138.572, 13.499, 9.951000000000001, 10298.25, 623.409, 1400, 9043.773999999999, 50.93, 914.462),
    @("FY3", 175.174, 91.06699999999999, 0.347, 0.005, 5968.532, -2997.343, 700, 5968.078, -2997.305, 700),
    @("总有效遮蔽并集时长(s)", 9.068)
)

$row = 11
foreach ($dataRow in $newRows) {
    $col = 1
    foreach ($val in $dataRow) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.Value = $val
        # Match the plain "vertical-center, no wrap" style already used by
        # the existing data rows (7-9) instead of inheriting the column's
        # default style.
        $cell.VerticalAlignment = -4108
        $col = $col + 1
    }
    $row = $row + 1
}
